$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.004031
$ws.Range("M2").Value = 3.087329333333333
$ws.Range("N2").Value = 9.261987999999999
$ws.Range("O2").Value = 0.1539049749041678
$ws.Range("P2").Value = 0.1539049749041678
$ws.Range("Q2").Value = 0.01244502454266666
$ws.Range("R2").Value = 0.112005220884
$ws.Range("S2").Value = 0.1539049749041678
$ws.Range("T2").Value = 0.1539049749041678

$ws.Range("G3").Value = 0.004031
$ws.Range("O3").Value = 0.2832552948356705
$ws.Range("P3").Value = 0.2832552948356705
$ws.Range("S3").Value = 0.2832552948356705
$ws.Range("T3").Value = 0.2832552948356705

$ws.Range("G4").Value = 0.004031
$ws.Range("M4").Value = 4.823431
$ws.Range("N4").Value = 14.470293
$ws.Range("O4").Value = 0.2404505470122564
$ws.Range("P4").Value = 0.2404505470122564
$ws.Range("Q4").Value = 0.019443250361
$ws.Range("R4").Value = 0.174989253249
$ws.Range("S4").Value = 0.2404505470122564
$ws.Range("T4").Value = 0.2404505470122564

$ws.Range("G5").Value = 0.004031
$ws.Range("M5").Value = 6.467117666666667
$ws.Range("N5").Value = 19.401353
$ws.Range("O5").Value = 0.3223891832479054
$ws.Range("P5").Value = 0.3223891832479053
$ws.Range("Q5").Value = 0.02606895131433333
$ws.Range("R5").Value = 0.234620561829
$ws.Range("S5").Value = 0.3223891832479054
$ws.Range("T5").Value = 0.3223891832479053
